# Update "want-to-go" counts (column F) and minimum ticket price (column G)
# on both the "展览" and "全部类型" worksheets - they carry duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 3066
    $ws.Range("F4").Value = 192
    $ws.Range("F7").Value = 1684
    $ws.Range("F12").Value = 1380
    $ws.Range("F14").Value = 524
    $ws.Range("F16").Value = 37
    $ws.Range("F17").Value = 7
    $ws.Range("F23").Value = 3233
    $ws.Range("F25").Value = 143
    $ws.Range("F26").Value = 329

    # Row 27 now has an actual minimum ticket price instead of "不可售"
    $ws.Range("F27").Value = 11
    $ws.Range("G27").Value = 50

    $ws.Range("F29").Value = 98
}
